$d = $word.ActiveDocument

$d.Content.Find.Execute("445÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "674÷8=", 2) | Out-Null
$d.Content.Find.Execute("779÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "701÷5=", 2) | Out-Null
$d.Content.Find.Execute("132÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "946÷7=", 2) | Out-Null
$d.Content.Find.Execute("798÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "194÷2=", 2) | Out-Null
$d.Content.Find.Execute("640÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "830÷3=", 2) | Out-Null
$d.Content.Find.Execute("297÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "495÷6=", 2) | Out-Null
$d.Content.Find.Execute("523÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "437÷5=", 2) | Out-Null
$d.Content.Find.Execute("363÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "522÷8=", 2) | Out-Null
$d.Content.Find.Execute("303÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "934÷8=", 2) | Out-Null
$d.Content.Find.Execute("591÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "509÷8=", 2) | Out-Null
$d.Content.Find.Execute("457÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "297÷8=", 2) | Out-Null
$d.Content.Find.Execute("116÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "614÷5=", 2) | Out-Null
$d.Content.Find.Execute("696÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "234÷4=", 2) | Out-Null
$d.Content.Find.Execute("432÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "447÷7=", 2) | Out-Null
$d.Content.Find.Execute("883÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "511÷2=", 2) | Out-Null
$d.Content.Find.Execute("430÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "600÷6=", 2) | Out-Null
$d.Content.Find.Execute("870÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "246÷4=", 2) | Out-Null
$d.Content.Find.Execute("944÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "239÷4=", 2) | Out-Null
$d.Content.Find.Execute("449÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "128÷2=", 2) | Out-Null
$d.Content.Find.Execute("439÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "712÷2=", 2) | Out-Null
$d.Content.Find.Execute("252÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "356÷8=", 2) | Out-Null
$d.Content.Find.Execute("165÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "848÷7=", 2) | Out-Null
$d.Content.Find.Execute("415÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "393÷2=", 2) | Out-Null
$d.Content.Find.Execute("305÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "936÷3=", 2) | Out-Null
$d.Content.Find.Execute("759÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "507÷4=", 2) | Out-Null
